$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.992.65"
$ws.Range("E2").Value = "  +0.57%  "

$ws.Range("D3").Value = "3.229.03"
$ws.Range("E3").Value = "  +3.02%  "

$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.18"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.63"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("D8").Value = "3.137.38"
$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  -1.34%  "

$ws.Range("E10").Value = "  -1.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.82"
$ws.Range("E11").Value = "  +2.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  -2.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  -2.87%  "

$ws.Range("E14").Value = "  +3.68%  "

$ws.Range("D15").Value = "3.745.95"
$ws.Range("E15").Value = "  +2.48%  "

$ws.Range("E16").Value = "  -1.33%  "

$ws.Range("D17").Value = "3.184.32"
$ws.Range("E17").Value = "  +1.60%  "

$ws.Range("D18").Value = "63.914.30"
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.06"
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "464.69"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.27"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.731"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.41"
$ws.Range("E23").Value = "  -1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.92"
$ws.Range("E24").Value = "  -2.90%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.00"
$ws.Range("E25").Value = "  -1.37%  "

$ws.Range("E26").Value = "  +3.78%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.30"
$ws.Range("E28").Value = "  +8.35%  "

$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").Value = "  +0.57%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.67"
$ws.Range("E30").Value = "  -0.98%  "

$ws.Range("E31").Value = "  +0.51%  "

$ws.Range("E32").Value = "  +4.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.03"
$ws.Range("E33").Value = "  +0.39%  "

$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("D35").Value = "0.0₃0860"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.32"
$ws.Range("E37").Value = "  -3.14%  "

$ws.Range("E38").Value = "  -1.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.99"
$ws.Range("E39").Value = "  -2.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "51.25"
$ws.Range("E40").Value = "  +0.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "437.04"
$ws.Range("E41").Value = "  -3.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.83"
$ws.Range("E42").Value = "  +1.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0371"
$ws.Range("E43").Value = "  -0.83%  "

$ws.Range("D44").Value = "2.911.23"
$ws.Range("E44").Value = "  -0.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.283"
$ws.Range("E45").Value = "  +1.55%  "

$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.87"
$ws.Range("E46").Value = "  +15.03%  "

$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.107"
$ws.Range("E47").Value = "  -3.72%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.41"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.19"
$ws.Range("E51").Value = "  +0.55%  "
